$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (aggregate/summary row) -----------------------------------------
$ws.Range("B2").Value = 1

$ws.Range("G2").Value = -24.83870967741936
$ws.Range("H2").Value = -24.83870967741936
$ws.Range("I2").Value = -16.7741935483871
$ws.Range("J2").Value = -16.7741935483871
$ws.Range("K2").Value = -2.48
$ws.Range("L2").Value = -20

$ws.Range("U2").Value = 14.8
$ws.Range("V2").Value = 0.1482965931863728
$ws.Range("W2").Value = -0.1530864197530864
$ws.Range("X2").Value = 0.09797387737246943
$ws.Range("Y2").Value = -0.2510602971255559
$ws.Range("Z2").Value = 0.01154562383612663
$ws.Range("AA2").Value = -0.1936685288640596
$ws.Range("AB2").Value = 0.08252131873405189
$ws.Range("AC2").Value = -0.2761898475981115
$ws.Range("AD2").Value = 34.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 34.4
$ws.Range("AG2").Value = 19.6
$ws.Range("AH2").Value = 0.2563338301043219
$ws.Range("AI2").Value = 0.7257383966244726
$ws.Range("AJ2").Value = 0.1641541038525963
$ws.Range("AK2").Value = 0.6012269938650308
$ws.Range("AL2").Value = 0.295
$ws.Range("AM2").Value = 0.09999999999999998
$ws.Range("AN2").Value = -18.59459459459459
$ws.Range("AO2").Value = -7.050847457627119
$ws.Range("AP2").Value = -10.59459459459459
$ws.Range("AQ2").Value = -20.8

# --- Row 4 (Renergen) is removed; its data moves up into row 3 -------------
$ws.Rows.Item(3).Delete()

# --- Row 3 now holds Renergen; refresh its name + figures -------------------
$ws.Range("B3").Value = "Renergen Limited (JSE:REN)"

$ws.Range("G3").Value = -24.83870967741936
$ws.Range("H3").Value = -24.83870967741936
$ws.Range("I3").Value = -16.7741935483871
$ws.Range("J3").Value = -16.7741935483871
$ws.Range("K3").Value = -2.48
$ws.Range("L3").Value = -20

$ws.Range("U3").Value = 14.8
$ws.Range("V3").Value = 0.1482965931863728
$ws.Range("W3").Value = -0.1530864197530864
$ws.Range("X3").Value = 0.09797387737246943
$ws.Range("Y3").Value = -0.2510602971255559
$ws.Range("Z3").Value = 0.01154562383612663
$ws.Range("AA3").Value = -0.1936685288640596
$ws.Range("AB3").Value = 0.08252131873405189
$ws.Range("AC3").Value = -0.2761898475981115
$ws.Range("AD3").Value = 34.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 34.4
$ws.Range("AG3").Value = 19.6
$ws.Range("AH3").Value = 0.2563338301043219
$ws.Range("AI3").Value = 0.7257383966244726
$ws.Range("AJ3").Value = 0.1641541038525963
$ws.Range("AK3").Value = 0.6012269938650308
$ws.Range("AL3").Value = 0.295
$ws.Range("AM3").Value = 0.09999999999999998
$ws.Range("AN3").Value = -18.59459459459459
$ws.Range("AO3").Value = -7.050847457627119
$ws.Range("AP3").Value = -10.59459459459459
$ws.Range("AQ3").Value = -20.8
